$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename shared string "Assessment" -> "Program evaluation" (row 12 label)
$ws.Range("A12").Value = "Program evaluation"

# Update the Publication/Assessment row values (row 12)
$ws.Range("C12").Value = 500
$ws.Range("D12").Value = 500
$ws.Range("E12").Value = 2000

# Update the selected cell to A12
$ws.Range("A12").Select()
